$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.345.03'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '1.868.55'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.77'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4696'
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2883'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06572'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.61'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07889'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.68'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').Value = '1.868.33'
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6923'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.107'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '267.85'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '30.282.67'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.98'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007696'
$ws.Range('E19').Value = '  +3.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.0000'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '2.109.94'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.235'
$ws.Range('E23').Value = '  -1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.189'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.408'
$ws.Range('E25').Value = '  +2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.68'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.360'
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09899'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.419'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.071'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04758'
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.134'
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7035'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.725'
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01875'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.799'
$ws.Range('E39').Value = '  +6.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.259'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.46'
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('E42').Value = '  +0.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4174'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('E45').Value = '  +0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.00'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '970.58'
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.117'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.080'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.51'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05677'
$ws.Range('E51').Value = '  +0.26%  '
